$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.165.92'
$ws.Range('E2').Value = '  +1.94%  '
$ws.Range('D3').Value = '3.440.58'
$ws.Range('E3').Value = '  +2.16%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '408.99'
$ws.Range('E5').Value = '  +0.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.03'
$ws.Range('E6').Value = '  -3.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.629'
$ws.Range('E7').Value = '  +6.56%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.750'
$ws.Range('E9').Value = '  +11.75%  '
$ws.Range('E10').Value = '  +18.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '43.41'
$ws.Range('E11').Value = '  +2.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.141'
$ws.Range('E12').Value = '  -0.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.78'
$ws.Range('E13').Value = '  +5.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.29'
$ws.Range('E14').Value = '  +3.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000198'
$ws.Range('E15').Value = '  +55.97%  '
$ws.Range('D16').Value = '3.422.32'
$ws.Range('E16').Value = '  +1.85%  '
$ws.Range('B17').Value = 'Polygon'
$ws.Range('C17').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.05'
$ws.Range('E17').Value = '  +3.08%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '62.135.39'
$ws.Range('E18').Value = '  +1.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.45'
$ws.Range('E19').Value = '  +3.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '374.32'
$ws.Range('E20').Value = '  +22.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '87.31'
$ws.Range('E21').Value = '  +4.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.19'
$ws.Range('E22').Value = '  -0.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.30'
$ws.Range('E23').Value = '  +4.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.22'
$ws.Range('E24').Value = '  +2.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '31.73'
$ws.Range('E25').Value = '  +8.00%  '
$ws.Range('E26').Value = '  +0.85%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.42'
$ws.Range('E27').Value = '  +1.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.71'
$ws.Range('E28').Value = '  +2.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.76'
$ws.Range('E29').Value = '  +11.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '44.11'
$ws.Range('E30').Value = '  +7.02%  '
$ws.Range('E31').Value = '  -0.73%  '
$ws.Range('E32').Value = '  -0.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.82'
$ws.Range('E33').Value = '  +4.57%  '
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('E35').Value = '  +2.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '52.13'
$ws.Range('E36').Value = '  +0.82%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.37'
$ws.Range('E38').Value = '  -1.65%  '
$ws.Range('E39').Value = '  +0.56%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.131'
$ws.Range('E40').Value = '  +6.66%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '143.36'
$ws.Range('E41').Value = '  +4.40%  '
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.313'
$ws.Range('E42').Value = '  +9.17%  '
$ws.Range('E43').Value = '  -0.38%  '
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.70'
$ws.Range('E45').Value = '  +0.42%  '
$ws.Range('E46').Value = '  +4.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '21.77'
$ws.Range('E47').Value = '  +1.70%  '
$ws.Range('D48').Value = '2.114.96'
$ws.Range('E48').Value = '  -0.24%  '
$ws.Range('E50').Value = '  +3.30%  '
$ws.Range('E51').Value = '  +7.00%  '
